# Projekt BD.xlsx - "Poprawki i porzadkowanie kodu."
#
# The G-column lookup list on the active sheet (G2, G4:G12) held:
#   Operator / Nazwa_operatora / Dzial / Uprawnienie / Data_konca_upr /
#   Rok / Mc / Dzien / Op_nazwisko / Op_imie
#
# "Nazwa_operatora" (G4) was a leftover/duplicate entry. It gets removed,
# and everything below it (G5:G12) shifts up one row, leaving G12 empty -
# the classic "select cell, Delete, Shift cells up" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift G5:G12 up into G4:G11 (content + formatting travel together),
# then clear the now-vacated last cell.
$ws.Range("G5:G12").Copy($ws.Range("G4")) | Out-Null
$ws.Range("G12").Clear() | Out-Null

# Match the resulting selection left behind by that edit.
$ws.Range("G4:G11").Select() | Out-Null
